$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on price cells whose new values would otherwise be auto-parsed as numbers
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D51').NumberFormat = "@"

# Apply updated values
$ws.Range('D2').Value = '26.851.16'
$ws.Range('E2').Value = '  -1.06%  '
$ws.Range('D3').Value = '1.563.20'
$ws.Range('E3').Value = '  +0.09%  '
$ws.Range('E4').Value = '  -0.11%  '
$ws.Range('D5').Value = '206.01'
$ws.Range('E5').Value = '  +0.03%  '
$ws.Range('D6').Value = '0.489'
$ws.Range('E6').Value = '  -0.78%  '
$ws.Range('E7').Value = '  -0.06%  '
$ws.Range('D8').Value = '21.74'
$ws.Range('E8').Value = '  -1.96%  '
$ws.Range('E9').Value = '  -0.08%  '
$ws.Range('E10').Value = '  -1.11%  '
$ws.Range('E11').Value = '  +0.54%  '
$ws.Range('D12').Value = '1.784.78'
$ws.Range('D13').Value = '1.567.53'
$ws.Range('E13').Value = '  +0.38%  '
$ws.Range('E14').Value = '  -1.01%  '
$ws.Range('D15').Value = '0.515'
$ws.Range('E15').Value = '  +0.07%  '
$ws.Range('D16').Value = '26.864.22'
$ws.Range('E16').Value = '  -0.94%  '
$ws.Range('D17').Value = '61.26'
$ws.Range('E17').Value = '  -2.58%  '
$ws.Range('D18').Value = '215.30'
$ws.Range('E18').Value = '  +1.48%  '
$ws.Range('E19').Value = '  +2.00%  '
$ws.Range('D20').Value = '0.0₃0680'
$ws.Range('E20').Value = '  -1.07%  '
$ws.Range('E21').Value = '  -0.10%  '
$ws.Range('E22').Value = '  +0.55%  '
$ws.Range('E23').Value = '  -2.20%  '
$ws.Range('E24').Value = '  +1.37%  '
$ws.Range('D25').Value = '153.88'
$ws.Range('E25').Value = '  +1.36%  '
$ws.Range('D26').Value = '6.73'
$ws.Range('E26').Value = '  +2.67%  '
$ws.Range('E28').Value = '  -0.07%  '
$ws.Range('E29').Value = '  -0.97%  '
$ws.Range('D30').Value = '0.0465'
$ws.Range('E30').Value = '  +0.41%  '
$ws.Range('E31').Value = '  -3.34%  '
$ws.Range('E32').Value = '  -0.01%  '
$ws.Range('D33').Value = '1.404.62'
$ws.Range('E33').Value = '  +2.15%  '
$ws.Range('E34').Value = '  -0.55%  '
$ws.Range('E35').Value = '  -1.35%  '
$ws.Range('E36').Value = '  -0.46%  '
$ws.Range('D37').Value = '0.916'
$ws.Range('E38').Value = '  -0.19%  '
$ws.Range('D39').Value = '0.528'
$ws.Range('E39').Value = '  +2.04%  '
$ws.Range('E41').Value = '  -0.08%  '
$ws.Range('E42').Value = '  +0.47%  '
$ws.Range('D43').Value = '5.45'
$ws.Range('E43').Value = '  +4.57%  '
$ws.Range('E44').Value = '  +0.31%  '
$ws.Range('B45').Value = 'Aave'
$ws.Range('C45').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D45').Value = '63.32'
$ws.Range('E45').Value = '  +0.04%  '
$ws.Range('B46').Value = 'RenderToken'
$ws.Range('C46').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D46').Value = '1.75'
$ws.Range('E46').Value = '  -1.21%  '
$ws.Range('D47').Value = '1.698.46'
$ws.Range('E47').Value = '  +0.18%  '
$ws.Range('D48').Value = '86.29'
$ws.Range('E48').Value = '  +1.03%  '
$ws.Range('D49').Value = '0.0505'
$ws.Range('E49').Value = '  +2.86%  '
$ws.Range('D50').Value = '0.0₇0978'
$ws.Range('E50').Value = '  -1.46%  '
$ws.Range('D51').Value = '0.0949'
$ws.Range('E51').Value = '  +1.08%  '

# Restore default style on the forced-text cells (keeps text value, drops explicit format)
$ws.Range('D5').Style = "Normal"
$ws.Range('D6').Style = "Normal"
$ws.Range('D8').Style = "Normal"
$ws.Range('D15').Style = "Normal"
$ws.Range('D17').Style = "Normal"
$ws.Range('D18').Style = "Normal"
$ws.Range('D25').Style = "Normal"
$ws.Range('D26').Style = "Normal"
$ws.Range('D30').Style = "Normal"
$ws.Range('D37').Style = "Normal"
$ws.Range('D39').Style = "Normal"
$ws.Range('D43').Style = "Normal"
$ws.Range('D45').Style = "Normal"
$ws.Range('D46').Style = "Normal"
$ws.Range('D48').Style = "Normal"
$ws.Range('D49').Style = "Normal"
$ws.Range('D51').Style = "Normal"
